# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the literal text "6-29-2011-12" for every
# data row; correct it to the proper ISO-style date text "2012-06-29".
#
# NumberFormat is forced to Text ("@") before the value is written so
# that Excel keeps the corrected value as a literal string instead of
# silently re-interpreting it as a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = $ws.Range("BF2:BF31")
$dateCol.NumberFormat = "@"
$dateCol.Value = "2012-06-29"
